$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (bold/border/center) from the existing header style to the new columns
$ws.Range("B1").Copy()
$ws.Range("AF1:AP1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B1").Value = 'filename'
$ws.Range("C1").Value = 'rays_present'
$ws.Range("D1").Value = 'approx_split'
$ws.Range("E1").Value = 'Diad1_pos'
$ws.Range("F1").Value = 'Diad2_pos'
$ws.Range("G1").Value = 'HB1_pos'
$ws.Range("H1").Value = 'HB2_pos'
$ws.Range("I1").Value = 'C13_pos'
$ws.Range("J1").Value = 'Diad1_abs_prom'
$ws.Range("K1").Value = 'Diad2_abs_prom'
$ws.Range("L1").Value = 'HB1_abs_prom'
$ws.Range("M1").Value = 'HB2_abs_prom'
$ws.Range("N1").Value = 'C13_abs_prom'
$ws.Range("O1").Value = 'Mean_abs_HB_prom'
$ws.Range("P1").Value = 'Diad2_HB2_abs_prom_ratio'
$ws.Range("Q1").Value = 'Diad1_HB1_abs_prom_ratio'
$ws.Range("R1").Value = 'Diad1_rel_prom'
$ws.Range("S1").Value = 'Diad2_rel_prom'
$ws.Range("T1").Value = 'HB1_rel_prom'
$ws.Range("U1").Value = 'HB2_rel_prom'
$ws.Range("V1").Value = 'C13_rel_prom'
$ws.Range("W1").Value = 'Diad1_HB1_abs_prom_ratio'
$ws.Range("X1").Value = 'Diad2_HB2_abs_prom_ratio'
$ws.Range("Y1").Value = 'Diad1_HB1_Valley_prom'
$ws.Range("Z1").Value = 'Diad2_HB2_abs_prom_ratio'
$ws.Range("AA1").Value = 'Mean_Diad_HB_Valley_prom'
$ws.Range("AB1").Value = 'Mean_abs_HB_prom'
$ws.Range("AC1").Value = 'Diad1_prom/std_betweendiads'
$ws.Range("AD1").Value = 'Diad2_prom/std_betweendiads'
$ws.Range("AE1").Value = 'Av_Diad_prom/std_betweendiads'
$ws.Range("AF1").Value = 'C13_prom/HB2_prom'
$ws.Range("AG1").Value = 'Av_Diad_HB_prom_ratio'
$ws.Range("AH1").Value = 'Diad2_height'
$ws.Range("AI1").Value = 'HB2_height'
$ws.Range("AJ1").Value = 'C13_height'
$ws.Range("AK1").Value = 'Diad1_height'
$ws.Range("AL1").Value = 'HB1_height'
$ws.Range("AM1").Value = 'Diad1_Median_Bck'
$ws.Range("AN1").Value = 'Diad2_Median_Bck'
$ws.Range("AO1").Value = 'C13_HB2_abs_prom_ratio'
$ws.Range("AP1").Value = 'Diad2_HB2_Valley_prom'
